$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test data values (avoid thread overwrite for test data values)
$ws.Range("C9").Value = "Test2"
$ws.Range("F9").Value = "y"

# Update the active selection to C9
$ws.Range("C9").Select()
